{"js": "// Update the worksheet date and the 25 division problems/answers in place.\n// Every old text value in the document is unique, so a matchCase, whole-\n// document search-and-replace for each pair is safe and unambiguous.\nconst replacements = [\n  [\"2026-01-14 Wednesday\", \"2026-01-15 Thursday\"],\n  [\"237\u00f79=26, 3\", \"811\u00f74=202, 3\"],\n  [\"884\u00f76=147, 2\", \"159\u00f75=31, 4\"],\n  [\"588\u00f76=98, 0\", \"426\u00f72=213, 0\"],\n  [\"733\u00f73=244, 1\", \"527\u00f78=65, 7\"],\n  [\"684\u00f78=85, 4\", \"846\u00f76=141, 0\"],\n  [\"591\u00f79=65, 6\", \"914\u00f79=101, 5\"],\n  [\"475\u00f79=52, 7\", \"858\u00f72=429, 0\"],\n  [\"922\u00f77=131, 5\", \"585\u00f79=65, 0\"],\n  [\"842\u00f73=280, 2\", \"803\u00f75=160, 3\"],\n  [\"826\u00f74=206, 2\", \"532\u00f72=266, 0\"],\n  [\"572\u00f76=95, 2\", \"535\u00f74=133, 3\"],\n  [\"526\u00f74=131, 2\", \"163\u00f75=32, 3\"],\n  [\"667\u00f79=74, 1\", \"447\u00f73=149, 0\"],\n  [\"749\u00f76=124, 5\", \"812\u00f76=135, 2\"],\n  [\"249\u00f77=35, 4\", \"653\u00f79=72, 5\"],\n  [\"503\u00f78=62, 7\", \"145\u00f72=72, 1\"],\n  [\"128\u00f73=42, 2\", \"356\u00f79=39, 5\"],\n  [\"123\u00f79=13, 6\", \"710\u00f73=236, 2\"],\n  [\"684\u00f77=97, 5\", \"314\u00f75=62, 4\"],\n  [\"443\u00f73=147, 2\", \"378\u00f79=42, 0\"],\n  [\"131\u00f78=16, 3\", \"349\u00f79=38, 7\"],\n  [\"949\u00f76=158, 1\", \"891\u00f76=148, 3\"],\n  [\"261\u00f73=87, 0\", \"645\u00f76=107, 3\"],\n  [\"929\u00f76=154, 5\", \"934\u00f79=103, 7\"],\n  [\"777\u00f78=97, 1\", \"394\u00f74=98, 2\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date and the 25 division problems/answers in place.\n# Every old text value in the document is unique, so a case-sensitive\n# Find/Replace (ReplaceAll) for each pair is safe and unambiguous.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2026-01-14 Wednesday\", \"2026-01-15 Thursday\"),\n    @(\"237\u00f79=26, 3\", \"811\u00f74=202, 3\"),\n    @(\"884\u00f76=147, 2\", \"159\u00f75=31, 4\"),\n    @(\"588\u00f76=98, 0\", \"426\u00f72=213, 0\"),\n    @(\"733\u00f73=244, 1\", \"527\u00f78=65, 7\"),\n    @(\"684\u00f78=85, 4\", \"846\u00f76=141, 0\"),\n    @(\"591\u00f79=65, 6\", \"914\u00f79=101, 5\"),\n    @(\"475\u00f79=52, 7\", \"858\u00f72=429, 0\"),\n    @(\"922\u00f77=131, 5\", \"585\u00f79=65, 0\"),\n    @(\"842\u00f73=280, 2\", \"803\u00f75=160, 3\"),\n    @(\"826\u00f74=206, 2\", \"532\u00f72=266, 0\"),\n    @(\"572\u00f76=95, 2\", \"535\u00f74=133, 3\"),\n    @(\"526\u00f74=131, 2\", \"163\u00f75=32, 3\"),\n    @(\"667\u00f79=74, 1\", \"447\u00f73=149, 0\"),\n    @(\"749\u00f76=124, 5\", \"812\u00f76=135, 2\"),\n    @(\"249\u00f77=35, 4\", \"653\u00f79=72, 5\"),\n    @(\"503\u00f78=62, 7\", \"145\u00f72=72, 1\"),\n    @(\"128\u00f73=42, 2\", \"356\u00f79=39, 5\"),\n    @(\"123\u00f79=13, 6\", \"710\u00f73=236, 2\"),\n    @(\"684\u00f77=97, 5\", \"314\u00f75=62, 4\"),\n    @(\"443\u00f73=147, 2\", \"378\u00f79=42, 0\"),\n    @(\"131\u00f78=16, 3\", \"349\u00f79=38, 7\"),\n    @(\"949\u00f76=158, 1\", \"891\u00f76=148, 3\"),\n    @(\"261\u00f73=87, 0\", \"645\u00f76=107, 3\"),\n    @(\"929\u00f76=154, 5\", \"934\u00f79=103, 7\"),\n    @(\"777\u00f78=97, 1\", \"394\u00f74=98, 2\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n\nWrite-Output \"done\"\n"}
